# "Deletei a porra toda" - remove the entire Von Neumann biography
# paragraph (the "C- Leonan von Octavius..." list item) along with the
# blank paragraph that precedes it, leaving only the trailing empty
# list paragraph (and the section properties) intact.

$d = $word.ActiveDocument

for ($i = $d.Paragraphs.Count; $i -ge 1; $i--) {
    $p = $d.Paragraphs($i)
    $t = $p.Range.Text

    # The paragraph holding the biography text about "Leonan von Octavius"
    $isBioParagraph = ($t -like "*Leonan*von*Octavius*") -or ($t -like "*Neumann falece*")

    # The leading blank paragraph that sits right before the biography
    # item (but never the final, required paragraph of the document).
    $isLeadingBlank = ($t.Trim() -eq "") -and ($i -lt $d.Paragraphs.Count)

    if ($isBioParagraph -or $isLeadingBlank) {
        $p.Range.Delete()
    }
}
